$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Row 10 (rule "R20"): change the "From" value (column C) from 18 to 1
$ws.Range("C10").Value = 1
